$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Pearson logo (image2.png -> image1.png) lives in both footers.
foreach ($ftr in $sec.Footers) {
    if ($ftr.Exists) {
        foreach ($shp in $ftr.Range.InlineShapes) {
            $shp.Name = "image1.png"
        }
    }
}

# BTec logo (image1.jpg -> image2.jpg) lives in both headers.
foreach ($hdr in $sec.Headers) {
    if ($hdr.Exists) {
        foreach ($shp in $hdr.Range.InlineShapes) {
            $shp.Name = "image2.jpg"
        }
    }
}
